$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch cell far outside the used range to build the literal date
# string as text (via a formula) so Excel's normal "looks like a date"
# auto-conversion does not turn it into a date serial number once it is
# copied into the real target cell.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""09/24/2025"""
$scratch.Copy()
$ws.Range("A23").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Range("B23").Value = 0.1319396418060274
$ws.Range("C23").Value = 0.8680603581939726
